$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the default-dictionary placeholder values that were parsed into row 6
# (F6:J6) — these should be empty rather than carrying the stale sample data.
$ws.Range("F6:J6").ClearContents()

# Update the active selection to match the new state.
$ws.Range("M17").Select()
